$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.616.80'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '2.442.84'
$ws.Range('E3').Value = '  -0.80%  '
$ws.Range('D5').Value = "'574.83"
$ws.Range('E5').Value = '  +0.77%  '
$ws.Range('D6').Value = "'144.48"
$ws.Range('E6').Value = '  -1.66%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('D9').Value = '2.439.23'
$ws.Range('E10').Value = '  -3.21%  '
$ws.Range('D11').Value = "'0.157"
$ws.Range('E11').Value = '  +0.48%  '
$ws.Range('D12').Value = "'5.22"
$ws.Range('E12').Value = '  -0.55%  '
$ws.Range('E13').Value = '  -1.32%  '
$ws.Range('D14').Value = "'26.71"
$ws.Range('E14').Value = '  -0.74%  '
$ws.Range('E15').Value = '  -2.05%  '
$ws.Range('D16').Value = '2.882.54'
$ws.Range('D17').Value = '62.365.53'
$ws.Range('E17').Value = '  -1.15%  '
$ws.Range('D18').Value = '2.440.32'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('D19').Value = "'11.16"
$ws.Range('E19').Value = '  -1.99%  '
$ws.Range('E20').Value = '  -2.00%  '
$ws.Range('D21').Value = "'328.93"
$ws.Range('D22').Value = "'4.16"
$ws.Range('E22').Value = '  -0.40%  '
$ws.Range('E23').Value = '  +4.34%  '
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').Value = "'65.57"
$ws.Range('E25').Value = '  -1.10%  '
$ws.Range('D26').Value = "'637.08"
$ws.Range('E26').Value = '  +2.08%  '
$ws.Range('D27').Value = "'9.22"
$ws.Range('E27').Value = '  +7.05%  '
$ws.Range('E28').Value = '  -5.46%  '
$ws.Range('D29').Value = '2.560.58'
$ws.Range('E29').Value = '  -0.28%  '
$ws.Range('D30').Value = "'0.999"
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('E31').Value = '  -2.99%  '
$ws.Range('D32').Value = "'8.10"
$ws.Range('E32').Value = '  -1.93%  '
$ws.Range('D33').Value = "'1.89"
$ws.Range('E33').Value = '  +0.57%  '
$ws.Range('E34').Value = '  -3.72%  '
$ws.Range('D35').Value = "'5.02"
$ws.Range('E35').Value = '  -1.12%  '
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('E37').Value = '  -1.65%  '
$ws.Range('E38').Value = '  -1.77%  '
$ws.Range('D39').Value = "'18.54"
$ws.Range('D40').Value = "'5.27"
$ws.Range('E40').Value = '  -2.76%  '
$ws.Range('D41').Value = "'146.15"
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('E42').Value = '  -2.50%  '
$ws.Range('D43').Value = "'42.30"
$ws.Range('E43').Value = '  +1.58%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = "'2.51"
$ws.Range('E45').Value = '  -3.84%  '
$ws.Range('D46').Value = "'145.75"
$ws.Range('E46').Value = '  -1.46%  '
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('E48').Value = '  -2.32%  '
$ws.Range('D49').Value = "'0.599"
$ws.Range('E49').Value = '  -0.53%  '
$ws.Range('D50').Value = "'19.71"
$ws.Range('E50').Value = '  -5.17%  '
$ws.Range('D51').Value = "'0.0231"
$ws.Range('E51').Value = '  -1.53%  '

# Reset style on cells that were set with a leading quote, to avoid leaving
# a text/quote-prefix style applied that was not in the original formatting.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
